# Add two new columns "I0" (I) and "IF" (J) to the sheet, with per-row
# numeric values, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1): copy the style from the existing header cell H1
#     so the new header cells look the same (bold, bordered, centered). ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2..69: set values for columns I and J ---
$iValues = @(8,8,7,9,10,6,7,9,6,8,7,6,9,5,3,7,5,8,7,6,7,7,7,5,8,10,6,7,5,8,9,9,7,6,6,5,7,8,6,10,7,6,9,6,6,11,7,7,6,5,3,6,9,6,6,5,9,9,7,4,5,9,8,8,8,5,4,5)
$jValues = @(8,8,8,9,10,6,7,9,6,9,8,7,9,5,4,7,7,9,7,7,7,7,7,6,8,10,6,7,6,8,9,9,7,6,8,7,7,8,6,11,7,6,9,7,7,11,8,7,6,5,5,6,9,7,6,7,9,9,8,4,6,9,8,8,8,7,6,6)

$startRow = 2
for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $startRow + $idx
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

$wb.Save()
